# Update "想去人数" (number of people interested) counts in the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.
#
# Each worksheet has the same events (by row) with the same F-column
# (想去人数 / interested-count) values that need bumping by a small
# amount, reflecting a newer data pull.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for the "展览" sheet (sheet name 展览)
$exhibitionUpdates = @{
    2  = 20610
    7  = 7713
    8  = 535
    12 = 175
    18 = 475
    19 = 76
    30 = 581
    31 = 101
    32 = 4921
    34 = 96
    36 = 12834
    43 = 4026
}

# Map of row number -> new value for the "全部类型" sheet
$allTypesUpdates = @{
    2  = 20610
    7  = 7713
    8  = 535
    12 = 175
    18 = 475
    19 = 76
    31 = 581
    33 = 101
    35 = 4921
    37 = 96
    39 = 12834
    46 = 4026
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
